$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to store the value as literal text (matching the
    # worksheet's existing inline-string / shared-string cells) instead of
    # letting Excel auto-convert numeric-looking strings (e.g. "1.007")
    # into a floating point number. Restore the original style afterwards
    # so no visible formatting/style change is introduced.
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '26.277.07'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '1.678.16'
$ws.Range('E3').Value = '  +0.61%  '
Set-TextValue 'D4' '1.007'
$ws.Range('E4').Value = '  +0.25%  '
Set-TextValue 'D5' '217.57'
$ws.Range('E5').Value = '  +0.28%  '
Set-TextValue 'D6' '0.5337'
$ws.Range('E6').Value = '  +4.43%  '
Set-TextValue 'D7' '1.006'
$ws.Range('E7').Value = '  +0.20%  '
Set-TextValue 'D8' '0.2680'
Set-TextValue 'D9' '0.06471'
$ws.Range('E9').Value = '  +0.99%  '
Set-TextValue 'D10' '21.87'
$ws.Range('E10').Value = '  -0.25%  '
Set-TextValue 'D11' '0.07544'
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').Value = '1.678.57'
$ws.Range('E12').Value = '  +0.62%  '
Set-TextValue 'D13' '4.517'
$ws.Range('E13').Value = '  +0.36%  '
Set-TextValue 'D14' '0.5764'
$ws.Range('E14').Value = '  -1.74%  '
Set-TextValue 'D15' '0.000008460'
$ws.Range('E15').Value = '  -1.07%  '
Set-TextValue 'D16' '64.64'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').Value = '26.281.41'
$ws.Range('E17').Value = '  +1.18%  '
Set-TextValue 'D18' '4.900'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  +0.22%  '
Set-TextValue 'D20' '10.86'
$ws.Range('E20').Value = '  +0.83%  '
Set-TextValue 'D21' '190.03'
$ws.Range('E21').Value = '  -0.35%  '
Set-TextValue 'D22' '6.198'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('E23').Value = '  +0.19%  '
Set-TextValue 'D24' '145.52'
$ws.Range('E24').Value = '  +0.13%  '
Set-TextValue 'D25' '7.810'
$ws.Range('E25').Value = '  +2.50%  '
Set-TextValue 'D26' '0.1267'
$ws.Range('E26').Value = '  +5.57%  '
Set-TextValue 'D27' '15.73'
$ws.Range('E27').Value = '  +0.61%  '
Set-TextValue 'D28' '0.06471'
$ws.Range('E28').Value = '  -4.07%  '
Set-TextValue 'D29' '1.383'
$ws.Range('E29').Value = '  +5.00%  '
Set-TextValue 'D30' '1.319'
$ws.Range('E30').Value = '  +0.15%  '
Set-TextValue 'D31' '3.583'
$ws.Range('E31').Value = '  +1.23%  '
Set-TextValue 'D32' '3.587'
$ws.Range('E32').Value = '  +1.82%  '
Set-TextValue 'D33' '1.656'
$ws.Range('E33').Value = '  +0.50%  '
Set-TextValue 'D34' '1.030'
$ws.Range('E34').Value = '  +1.27%  '
Set-TextValue 'D35' '0.6184'
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('E36').Value = '  +1.54%  '
Set-TextValue 'D37' '2.720'
$ws.Range('E37').Value = '  +0.11%  '
Set-TextValue 'D38' '6.240'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '1.111.02'
$ws.Range('E39').Value = '  +2.21%  '
Set-TextValue 'D40' '0.01621'
$ws.Range('E40').Value = '  +1.09%  '
Set-TextValue 'D41' '0.8715'
$ws.Range('E41').Value = '  +0.28%  '
Set-TextValue 'D42' '1.013'
$ws.Range('E42').Value = '  +0.40%  '
Set-TextValue 'D43' '100.30'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('D44').Value = '1.828.59'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('E45').Value = '  -3.52%  '
Set-TextValue 'D46' '57.06'
$ws.Range('E46').Value = '  +1.18%  '
Set-TextValue 'D47' '8.138'
$ws.Range('E47').Value = '  +0.54%  '
Set-TextValue 'D48' '1.004'
$ws.Range('E48').Value = '  -0.15%  '
Set-TextValue 'D49' '0.05263'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('E50').Value = '  +1.19%  '
Set-TextValue 'D51' '0.4288'
